$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper cell used to stage text-formatted numeric-looking strings so that
# PasteSpecial(xlPasteValues) can drop them into the target cells as plain
# text without Excel's automatic text-to-number coercion, and without
# leaving a new NumberFormat style on the destination cell.
$helper = $ws.Range("ZZ1")
$helper.NumberFormat = "@"

$ws.Range("D2").Value = '69.001.61'
$ws.Range("E2").Value = '  +2.70%  '
$ws.Range("D3").Value = '3.735.99'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("E4").Value = '  +0.20%  '
$helper.Value = '601.23'
$helper.Copy() | Out-Null
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = '  +1.71%  '
$helper.Value = '168.64'
$helper.Copy() | Out-Null
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("D7").Value = '3.736.37'
$ws.Range("E7").Value = '  -0.83%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("E9").Value = '  +2.40%  '
$helper.Value = '0.164'
$helper.Copy() | Out-Null
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = '  +4.03%  '
$helper.Value = '6.27'
$helper.Copy() | Out-Null
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = '  -0.18%  '
$ws.Range("E12").Value = '  +0.39%  '
$helper.Value = '38.19'
$helper.Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = '  +1.29%  '
$ws.Range("E14").Value = '  +0.40%  '
$ws.Range("D15").Value = '4.361.35'
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '3.742.51'
$ws.Range("E16").Value = '  -0.10%  '
$ws.Range("D17").Value = '68.934.83'
$ws.Range("E17").Value = '  +2.47%  '
$helper.Value = '7.25'
$helper.Copy() | Out-Null
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = '  +2.44%  '
$ws.Range("E19").Value = '  +0.31%  '
$helper.Value = '17.20'
$helper.Copy() | Out-Null
$ws.Range("D20").PasteSpecial(-4163)
$ws.Range("E20").Value = '  +6.96%  '
$helper.Value = '497.28'
$helper.Copy() | Out-Null
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = '  +2.33%  '
$helper.Value = '9.51'
$helper.Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = '  +4.63%  '
$helper.Value = '0.723'
$helper.Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = '  +0.39%  '
$helper.Value = '84.76'
$helper.Copy() | Out-Null
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = '  +0.97%  '
$helper.Value = '0.0000142'
$helper.Copy() | Out-Null
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = '  +3.76%  '
$helper.Value = '2.31'
$helper.Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4163)
$ws.Range("E26").Value = '  -1.78%  '
$helper.Value = '12.28'
$helper.Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = '  +0.69%  '
$helper.Value = '10.11'
$helper.Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("E30").Value = '  +0.88%  '
$helper.Value = '2.42'
$helper.Copy() | Out-Null
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = '  +0.96%  '
$helper.Value = '7.95'
$helper.Copy() | Out-Null
$ws.Range("D32").PasteSpecial(-4163)
$ws.Range("E32").Value = '  +3.22%  '
$helper.Value = '31.62'
$helper.Copy() | Out-Null
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = '  -1.65%  '
$ws.Range("D34").Value = '3.877.50'
$ws.Range("E34").Value = '  -0.08%  '
$ws.Range("E35").Value = '  +1.52%  '
$ws.Range("D36").Value = '3.677.25'
$ws.Range("E36").Value = '  -0.44%  '
$ws.Range("E37").Value = '  +0.16%  '
$ws.Range("E38").Value = '  +0.86%  '
$helper.Value = '5.79'
$helper.Copy() | Out-Null
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = '  +1.32%  '
$helper.Value = '0.133'
$helper.Copy() | Out-Null
$ws.Range("D40").PasteSpecial(-4163)
$ws.Range("E40").Value = '  -1.24%  '
$ws.Range("E41").Value = '  +0.25%  '
$helper.Value = '436.47'
$helper.Copy() | Out-Null
$ws.Range("D42").PasteSpecial(-4163)
$ws.Range("E42").Value = '  -3.42%  '
$helper.Value = '49.00'
$helper.Copy() | Out-Null
$ws.Range("D43").PasteSpecial(-4163)
$ws.Range("E43").Value = '  +0.52%  '
$helper.Value = '1.99'
$helper.Copy() | Out-Null
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = '  +1.26%  '
$ws.Range("E45").Value = '  +0.74%  '
$helper.Value = '8.38'
$helper.Copy() | Out-Null
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = '  +1.98%  '
$ws.Range("E47").Value = '  +0.00%  '
$helper.Value = '40.40'
$helper.Copy() | Out-Null
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = '  -1.98%  '
$ws.Range("E49").Value = '  +3.18%  '
$ws.Range("E50").Value = '  +1.50%  '
$ws.Range("D51").Value = '2.754.85'
$ws.Range("E51").Value = '  -1.76%  '

$helper.Clear()
$excel.CutCopyMode = $false
